$wb = $excel.ActiveWorkbook

# Sheet1 (大智投资组合收益) - refresh existing rows 11-20 (timestamp 202506092206)
$ws = $wb.Worksheets.Item(1)
$ws.Range("B11:B20").NumberFormat = "@"
$ws.Range("I11:I20").NumberFormat = "@"
$ws.Range("A11").Value = '大智 (稳健智远)'
$ws.Range("B11").Value = '000333'
$ws.Range("C11").Value = '美的集团'
$ws.Range("D11").Value = 3.0
$ws.Range("E11").Value = 40.00502976847048
$ws.Range("F11").Value = 75.05
$ws.Range("G11").Value = 3002.377484123709
$ws.Range("H11").Value = 100079.2494707903
$ws.Range("I11").Value = '202506092206'
$ws.Range("A12").Value = '大智 (稳健智远)'
$ws.Range("B12").Value = '510050'
$ws.Range("C12").Value = '上证50ETF'
$ws.Range("D12").Value = 5.0
$ws.Range("E12").Value = 1826.263676474276
$ws.Range("F12").Value = 2.74
$ws.Range("G12").Value = 5003.962473539516
$ws.Range("H12").Value = 100079.2494707903
$ws.Range("I12").Value = '202506092206'
$ws.Range("A13").Value = '大智 (稳健智远)'
$ws.Range("B13").Value = '510300'
$ws.Range("C13").Value = '沪深300ETF'
$ws.Range("D13").Value = 5.0
$ws.Range("E13").Value = 1260.44394799484
$ws.Range("F13").Value = 3.97
$ws.Range("G13").Value = 5003.962473539516
$ws.Range("H13").Value = 100079.2494707903
$ws.Range("I13").Value = '202506092206'
$ws.Range("A14").Value = '大智 (稳健智远)'
$ws.Range("B14").Value = '518880'
$ws.Range("C14").Value = '黄金ETF'
$ws.Range("D14").Value = 5.0
$ws.Range("E14").Value = 675.2985794250359
$ws.Range("F14").Value = 7.41
$ws.Range("G14").Value = 5003.962473539516
$ws.Range("H14").Value = 100079.2494707903
$ws.Range("I14").Value = '202506092206'
$ws.Range("A15").Value = '大智 (稳健智远)'
$ws.Range("B15").Value = '600085'
$ws.Range("C15").Value = '同仁堂'
$ws.Range("D15").Value = 2.0
$ws.Range("E15").Value = 53.02211892492203
$ws.Range("F15").Value = 37.75
$ws.Range("G15").Value = 2001.584989415806
$ws.Range("H15").Value = 100079.2494707903
$ws.Range("I15").Value = '202506092206'
$ws.Range("A16").Value = '大智 (稳健智远)'
$ws.Range("B16").Value = '600900'
$ws.Range("C16").Value = '长江电力'
$ws.Range("D16").Value = 20.0
$ws.Range("E16").Value = 666.5284680039316
$ws.Range("F16").Value = 30.03
$ws.Range("G16").Value = 20015.84989415807
$ws.Range("H16").Value = 100079.2494707903
$ws.Range("I16").Value = '202506092206'
$ws.Range("A17").Value = '大智 (稳健智远)'
$ws.Range("B17").Value = '600989'
$ws.Range("C17").Value = '宝丰能源'
$ws.Range("D17").Value = 5.0
$ws.Range("E17").Value = 310.9982892193609
$ws.Range("F17").Value = 16.09
$ws.Range("G17").Value = 5003.962473539516
$ws.Range("H17").Value = 100079.2494707903
$ws.Range("I17").Value = '202506092206'
$ws.Range("A18").Value = '大智 (稳健智远)'
$ws.Range("B18").Value = 'HK02899'
$ws.Range("C18").Value = '紫金矿业'
$ws.Range("D18").Value = 20.0
$ws.Range("E18").Value = 1087.81792903033
$ws.Range("F18").Value = 18.4
$ws.Range("G18").Value = 20015.84989415807
$ws.Range("H18").Value = 100079.2494707903
$ws.Range("I18").Value = '202506092206'
$ws.Range("A19").Value = '大智 (稳健智远)'
$ws.Range("B19").Value = 'HK06881'
$ws.Range("C19").Value = '中国银河'
$ws.Range("D19").Value = 5.0
$ws.Range("E19").Value = 616.2515361501868
$ws.Range("F19").Value = 8.12
$ws.Range("G19").Value = 5003.962473539516
$ws.Range("H19").Value = 100079.2494707903
$ws.Range("I19").Value = '202506092206'
$ws.Range("A20").Value = '大智 (稳健智远)'
$ws.Range("B20").Value = '100000'
$ws.Range("C20").Value = '现金'
$ws.Range("D20").Value = 30.0
$ws.Range("E20").Value = 30023.7748412371
$ws.Range("F20").Value = 1.0
$ws.Range("G20").Value = 30023.7748412371
$ws.Range("H20").Value = 100079.2494707903
$ws.Range("I20").Value = '202506092206'
$ws.Range("B11:B20").ClearFormats()
$ws.Range("I11:I20").ClearFormats()

# Sheet1 (大智投资组合收益) - append new rows 21-30 (timestamp 202506101600)
$ws = $wb.Worksheets.Item(1)
$ws.Range("B21:B30").NumberFormat = "@"
$ws.Range("I21:I30").NumberFormat = "@"
$ws.Range("A21").Value = '大智 (稳健智远)'
$ws.Range("B21").Value = '000333'
$ws.Range("C21").Value = '美的集团'
$ws.Range("D21").Value = 1.499999999999992
$ws.Range("E21").Value = 40.00502976847048
$ws.Range("F21").Value = 75.05
$ws.Range("G21").Value = 3002.377484123709
$ws.Range("H21").Value = 200158.4989415806
$ws.Range("I21").Value = '202506101600'
$ws.Range("A22").Value = '大智 (稳健智远)'
$ws.Range("B22").Value = '510050'
$ws.Range("C22").Value = '上证50ETF'
$ws.Range("D22").Value = 2.499999999999988
$ws.Range("E22").Value = 1826.263676474276
$ws.Range("F22").Value = 2.74
$ws.Range("G22").Value = 5003.962473539516
$ws.Range("H22").Value = 200158.4989415806
$ws.Range("I22").Value = '202506101600'
$ws.Range("A23").Value = '大智 (稳健智远)'
$ws.Range("B23").Value = '510300'
$ws.Range("C23").Value = '沪深300ETF'
$ws.Range("D23").Value = 2.499999999999988
$ws.Range("E23").Value = 1260.44394799484
$ws.Range("F23").Value = 3.97
$ws.Range("G23").Value = 5003.962473539516
$ws.Range("H23").Value = 200158.4989415806
$ws.Range("I23").Value = '202506101600'
$ws.Range("A24").Value = '大智 (稳健智远)'
$ws.Range("B24").Value = '518880'
$ws.Range("C24").Value = '黄金ETF'
$ws.Range("D24").Value = 2.499999999999988
$ws.Range("E24").Value = 675.2985794250359
$ws.Range("F24").Value = 7.41
$ws.Range("G24").Value = 5003.962473539516
$ws.Range("H24").Value = 200158.4989415806
$ws.Range("I24").Value = '202506101600'
$ws.Range("A25").Value = '大智 (稳健智远)'
$ws.Range("B25").Value = '600085'
$ws.Range("C25").Value = '同仁堂'
$ws.Range("D25").Value = 0.999999999999995
$ws.Range("E25").Value = 53.02211892492203
$ws.Range("F25").Value = 37.75
$ws.Range("G25").Value = 2001.584989415806
$ws.Range("H25").Value = 200158.4989415806
$ws.Range("I25").Value = '202506101600'
$ws.Range("A26").Value = '大智 (稳健智远)'
$ws.Range("B26").Value = '600900'
$ws.Range("C26").Value = '长江电力'
$ws.Range("D26").Value = 9.99999999999995
$ws.Range("E26").Value = 666.5284680039316
$ws.Range("F26").Value = 30.03
$ws.Range("G26").Value = 20015.84989415807
$ws.Range("H26").Value = 200158.4989415806
$ws.Range("I26").Value = '202506101600'
$ws.Range("A27").Value = '大智 (稳健智远)'
$ws.Range("B27").Value = '600989'
$ws.Range("C27").Value = '宝丰能源'
$ws.Range("D27").Value = 2.499999999999988
$ws.Range("E27").Value = 310.9982892193609
$ws.Range("F27").Value = 16.09
$ws.Range("G27").Value = 5003.962473539516
$ws.Range("H27").Value = 200158.4989415806
$ws.Range("I27").Value = '202506101600'
$ws.Range("A28").Value = '大智 (稳健智远)'
$ws.Range("B28").Value = 'HK02899'
$ws.Range("C28").Value = '紫金矿业'
$ws.Range("D28").Value = 9.99999999999995
$ws.Range("E28").Value = 1087.81792903033
$ws.Range("F28").Value = 18.4
$ws.Range("G28").Value = 20015.84989415807
$ws.Range("H28").Value = 200158.4989415806
$ws.Range("I28").Value = '202506101600'
$ws.Range("A29").Value = '大智 (稳健智远)'
$ws.Range("B29").Value = 'HK06881'
$ws.Range("C29").Value = '中国银河'
$ws.Range("D29").Value = 2.499999999999988
$ws.Range("E29").Value = 616.2515361501868
$ws.Range("F29").Value = 8.12
$ws.Range("G29").Value = 5003.962473539516
$ws.Range("H29").Value = 200158.4989415806
$ws.Range("I29").Value = '202506101600'
$ws.Range("A30").Value = '大智 (稳健智远)'
$ws.Range("B30").Value = '100000'
$ws.Range("C30").Value = '现金'
$ws.Range("D30").Value = 14.99999999999993
$ws.Range("E30").Value = 30023.7748412371
$ws.Range("F30").Value = 1.0
$ws.Range("G30").Value = 30023.7748412371
$ws.Range("H30").Value = 200158.4989415806
$ws.Range("I30").Value = '202506101600'
$ws.Range("B21:B30").ClearFormats()
$ws.Range("I21:I30").ClearFormats()

# Sheet2 (大成投资组合收益) - refresh existing rows 9-15 (timestamp 202506092206)
$ws = $wb.Worksheets.Item(2)
$ws.Range("B9:B15").NumberFormat = "@"
$ws.Range("I9:I15").NumberFormat = "@"
$ws.Range("A9").Value = '大成 (锐进先锋)'
$ws.Range("B9").Value = '000725'
$ws.Range("C9").Value = '京东方A'
$ws.Range("D9").Value = 5.0
$ws.Range("E9").Value = 1222.947139927285
$ws.Range("F9").Value = 3.9
$ws.Range("G9").Value = 4769.493845716412
$ws.Range("H9").Value = 95389.87691432823
$ws.Range("I9").Value = '202506092206'
$ws.Range("A10").Value = '大成 (锐进先锋)'
$ws.Range("B10").Value = '159781'
$ws.Range("C10").Value = '科创创业ETF'
$ws.Range("D10").Value = 5.0
$ws.Range("E10").Value = 8999.044991917757
$ws.Range("F10").Value = 0.53
$ws.Range("G10").Value = 4769.493845716412
$ws.Range("H10").Value = 95389.87691432823
$ws.Range("I10").Value = '202506092206'
$ws.Range("A11").Value = '大成 (锐进先锋)'
$ws.Range("B11").Value = '513100'
$ws.Range("C11").Value = '纳指ETF'
$ws.Range("D11").Value = 5.0
$ws.Range("E11").Value = 3037.894169246122
$ws.Range("F11").Value = 1.57
$ws.Range("G11").Value = 4769.493845716412
$ws.Range("H11").Value = 95389.87691432823
$ws.Range("I11").Value = '202506092206'
$ws.Range("A12").Value = '大成 (锐进先锋)'
$ws.Range("B12").Value = '513290'
$ws.Range("C12").Value = '纳指生物科技ETF'
$ws.Range("D12").Value = 1.0
$ws.Range("E12").Value = 851.6953295922162
$ws.Range("F12").Value = 1.12
$ws.Range("G12").Value = 953.8987691432823
$ws.Range("H12").Value = 95389.87691432823
$ws.Range("I12").Value = '202506092206'
$ws.Range("A13").Value = '大成 (锐进先锋)'
$ws.Range("B13").Value = '603119'
$ws.Range("C13").Value = '浙江荣泰'
$ws.Range("D13").Value = 45.0
$ws.Range("E13").Value = 1080.701022443296
$ws.Range("F13").Value = 39.72
$ws.Range("G13").Value = 42925.4446114477
$ws.Range("H13").Value = 95389.87691432823
$ws.Range("I13").Value = '202506092206'
$ws.Range("A14").Value = '大成 (锐进先锋)'
$ws.Range("B14").Value = '688290'
$ws.Range("C14").Value = '景业智能'
$ws.Range("D14").Value = 9.0
$ws.Range("E14").Value = 163.4632315744391
$ws.Range("F14").Value = 52.52
$ws.Range("G14").Value = 8585.08892228954
$ws.Range("H14").Value = 95389.87691432823
$ws.Range("I14").Value = '202506092206'
$ws.Range("A15").Value = '大成 (锐进先锋)'
$ws.Range("B15").Value = '100000'
$ws.Range("C15").Value = '现金'
$ws.Range("D15").Value = 30.0
$ws.Range("E15").Value = 28616.96307429847
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 28616.96307429847
$ws.Range("H15").Value = 95389.87691432823
$ws.Range("I15").Value = '202506092206'
$ws.Range("B9:B15").ClearFormats()
$ws.Range("I9:I15").ClearFormats()

# Sheet2 (大成投资组合收益) - append new rows 16-22 (timestamp 202506101600)
$ws = $wb.Worksheets.Item(2)
$ws.Range("B16:B22").NumberFormat = "@"
$ws.Range("I16:I22").NumberFormat = "@"
$ws.Range("A16").Value = '大成 (锐进先锋)'
$ws.Range("B16").Value = '000725'
$ws.Range("C16").Value = '京东方A'
$ws.Range("D16").Value = 2.499999999999987
$ws.Range("E16").Value = 1222.947139927285
$ws.Range("F16").Value = 3.9
$ws.Range("G16").Value = 4769.493845716412
$ws.Range("H16").Value = 190779.7538286565
$ws.Range("I16").Value = '202506101600'
$ws.Range("A17").Value = '大成 (锐进先锋)'
$ws.Range("B17").Value = '159781'
$ws.Range("C17").Value = '科创创业ETF'
$ws.Range("D17").Value = 2.499999999999987
$ws.Range("E17").Value = 8999.044991917757
$ws.Range("F17").Value = 0.53
$ws.Range("G17").Value = 4769.493845716412
$ws.Range("H17").Value = 190779.7538286565
$ws.Range("I17").Value = '202506101600'
$ws.Range("A18").Value = '大成 (锐进先锋)'
$ws.Range("B18").Value = '513100'
$ws.Range("C18").Value = '纳指ETF'
$ws.Range("D18").Value = 2.499999999999987
$ws.Range("E18").Value = 3037.894169246122
$ws.Range("F18").Value = 1.57
$ws.Range("G18").Value = 4769.493845716412
$ws.Range("H18").Value = 190779.7538286565
$ws.Range("I18").Value = '202506101600'
$ws.Range("A19").Value = '大成 (锐进先锋)'
$ws.Range("B19").Value = '513290'
$ws.Range("C19").Value = '纳指生物科技ETF'
$ws.Range("D19").Value = 0.4999999999999974
$ws.Range("E19").Value = 851.6953295922162
$ws.Range("F19").Value = 1.12
$ws.Range("G19").Value = 953.8987691432823
$ws.Range("H19").Value = 190779.7538286565
$ws.Range("I19").Value = '202506101600'
$ws.Range("A20").Value = '大成 (锐进先锋)'
$ws.Range("B20").Value = '603119'
$ws.Range("C20").Value = '浙江荣泰'
$ws.Range("D20").Value = 22.49999999999988
$ws.Range("E20").Value = 1080.701022443296
$ws.Range("F20").Value = 39.72
$ws.Range("G20").Value = 42925.4446114477
$ws.Range("H20").Value = 190779.7538286565
$ws.Range("I20").Value = '202506101600'
$ws.Range("A21").Value = '大成 (锐进先锋)'
$ws.Range("B21").Value = '688290'
$ws.Range("C21").Value = '景业智能'
$ws.Range("D21").Value = 4.499999999999977
$ws.Range("E21").Value = 163.4632315744391
$ws.Range("F21").Value = 52.52
$ws.Range("G21").Value = 8585.08892228954
$ws.Range("H21").Value = 190779.7538286565
$ws.Range("I21").Value = '202506101600'
$ws.Range("A22").Value = '大成 (锐进先锋)'
$ws.Range("B22").Value = '100000'
$ws.Range("C22").Value = '现金'
$ws.Range("D22").Value = 14.99999999999992
$ws.Range("E22").Value = 28616.96307429847
$ws.Range("F22").Value = 1.0
$ws.Range("G22").Value = 28616.96307429847
$ws.Range("H22").Value = 190779.7538286565
$ws.Range("I22").Value = '202506101600'
$ws.Range("B16:B22").ClearFormats()
$ws.Range("I16:I22").ClearFormats()

# Sheet3 (我的投资组合收益) - refresh existing rows 16-30 (timestamp 202506092206)
$ws = $wb.Worksheets.Item(3)
$ws.Range("B16:B30").NumberFormat = "@"
$ws.Range("I16:I30").NumberFormat = "@"
$ws.Range("A16").Value = '范式进化投资组合'
$ws.Range("B16").Value = '000333'
$ws.Range("C16").Value = '美的集团'
$ws.Range("D16").Value = 1.0
$ws.Range("E16").Value = 13.28200178359967
$ws.Range("F16").Value = 75.05
$ws.Range("G16").Value = 996.8142338591555
$ws.Range("H16").Value = 99681.42338591556
$ws.Range("I16").Value = '202506092206'
$ws.Range("A17").Value = '范式进化投资组合'
$ws.Range("B17").Value = '000725'
$ws.Range("C17").Value = '京东方A'
$ws.Range("D17").Value = 5.0
$ws.Range("E17").Value = 1277.966966486097
$ws.Range("F17").Value = 3.9
$ws.Range("G17").Value = 4984.071169295778
$ws.Range("H17").Value = 99681.42338591556
$ws.Range("I17").Value = '202506092206'
$ws.Range("A18").Value = '范式进化投资组合'
$ws.Range("B18").Value = '159781'
$ws.Range("C18").Value = '科创创业ETF'
$ws.Range("D18").Value = 5.0
$ws.Range("E18").Value = 9403.907866595806
$ws.Range("F18").Value = 0.53
$ws.Range("G18").Value = 4984.071169295778
$ws.Range("H18").Value = 99681.42338591556
$ws.Range("I18").Value = '202506092206'
$ws.Range("A19").Value = '范式进化投资组合'
$ws.Range("B19").Value = '510050'
$ws.Range("C19").Value = '上证50ETF'
$ws.Range("D19").Value = 5.0
$ws.Range("E19").Value = 1819.00407638532
$ws.Range("F19").Value = 2.74
$ws.Range("G19").Value = 4984.071169295778
$ws.Range("H19").Value = 99681.42338591556
$ws.Range("I19").Value = '202506092206'
$ws.Range("A20").Value = '范式进化投资组合'
$ws.Range("B20").Value = '510300'
$ws.Range("C20").Value = '沪深300ETF'
$ws.Range("D20").Value = 5.0
$ws.Range("E20").Value = 1255.433543903218
$ws.Range("F20").Value = 3.97
$ws.Range("G20").Value = 4984.071169295778
$ws.Range("H20").Value = 99681.42338591556
$ws.Range("I20").Value = '202506092206'
$ws.Range("A21").Value = '范式进化投资组合'
$ws.Range("B21").Value = '513100'
$ws.Range("C21").Value = '纳指ETF'
$ws.Range("D21").Value = 1.0
$ws.Range("E21").Value = 634.9135247510544
$ws.Range("F21").Value = 1.57
$ws.Range("G21").Value = 996.8142338591555
$ws.Range("H21").Value = 99681.42338591556
$ws.Range("I21").Value = '202506092206'
$ws.Range("A22").Value = '范式进化投资组合'
$ws.Range("B22").Value = '513290'
$ws.Range("C22").Value = '纳指生物科技ETF'
$ws.Range("D22").Value = 1.0
$ws.Range("E22").Value = 890.0127088028173
$ws.Range("F22").Value = 1.12
$ws.Range("G22").Value = 996.8142338591555
$ws.Range("H22").Value = 99681.42338591556
$ws.Range("I22").Value = '202506092206'
$ws.Range("A23").Value = '范式进化投资组合'
$ws.Range("B23").Value = '518880'
$ws.Range("C23").Value = '黄金ETF'
$ws.Range("D23").Value = 1.0
$ws.Range("E23").Value = 134.5228385774839
$ws.Range("F23").Value = 7.41
$ws.Range("G23").Value = 996.8142338591555
$ws.Range("H23").Value = 99681.42338591556
$ws.Range("I23").Value = '202506092206'
$ws.Range("A24").Value = '范式进化投资组合'
$ws.Range("B24").Value = '600085'
$ws.Range("C24").Value = '同仁堂'
$ws.Range("D24").Value = 1.0
$ws.Range("E24").Value = 26.4056750691167
$ws.Range("F24").Value = 37.75
$ws.Range("G24").Value = 996.8142338591555
$ws.Range("H24").Value = 99681.42338591556
$ws.Range("I24").Value = '202506092206'
$ws.Range("A25").Value = '范式进化投资组合'
$ws.Range("B25").Value = '600900'
$ws.Range("C25").Value = '长江电力'
$ws.Range("D25").Value = 1.0
$ws.Range("E25").Value = 33.19394718145706
$ws.Range("F25").Value = 30.03
$ws.Range("G25").Value = 996.8142338591554
$ws.Range("H25").Value = 99681.42338591556
$ws.Range("I25").Value = '202506092206'
$ws.Range("A26").Value = '范式进化投资组合'
$ws.Range("B26").Value = '600989'
$ws.Range("C26").Value = '宝丰能源'
$ws.Range("D26").Value = 5.0
$ws.Range("E26").Value = 309.7620366249706
$ws.Range("F26").Value = 16.09
$ws.Range("G26").Value = 4984.071169295778
$ws.Range("H26").Value = 99681.42338591556
$ws.Range("I26").Value = '202506092206'
$ws.Range("A27").Value = '范式进化投资组合'
$ws.Range("B27").Value = '603119'
$ws.Range("C27").Value = '浙江荣泰'
$ws.Range("D27").Value = 1.0
$ws.Range("E27").Value = 25.09602804277834
$ws.Range("F27").Value = 39.72
$ws.Range("G27").Value = 996.8142338591555
$ws.Range("H27").Value = 99681.42338591556
$ws.Range("I27").Value = '202506092206'
$ws.Range("A28").Value = '范式进化投资组合'
$ws.Range("B28").Value = 'HK02899'
$ws.Range("C28").Value = '紫金矿业'
$ws.Range("D28").Value = 1.0
$ws.Range("E28").Value = 54.17468662278019
$ws.Range("F28").Value = 18.4
$ws.Range("G28").Value = 996.8142338591555
$ws.Range("H28").Value = 99681.42338591556
$ws.Range("I28").Value = '202506092206'
$ws.Range("A29").Value = '范式进化投资组合'
$ws.Range("B29").Value = 'HK06881'
$ws.Range("C29").Value = '中国银河'
$ws.Range("D29").Value = 1.0
$ws.Range("E29").Value = 122.7603736279748
$ws.Range("F29").Value = 8.12
$ws.Range("G29").Value = 996.8142338591555
$ws.Range("H29").Value = 99681.42338591556
$ws.Range("I29").Value = '202506092206'
$ws.Range("A30").Value = '范式进化投资组合'
$ws.Range("B30").Value = '100000'
$ws.Range("C30").Value = '现金'
$ws.Range("D30").Value = 66.0
$ws.Range("E30").Value = 65789.73943470427
$ws.Range("F30").Value = 1.0
$ws.Range("G30").Value = 65789.73943470427
$ws.Range("H30").Value = 99681.42338591556
$ws.Range("I30").Value = '202506092206'
$ws.Range("B16:B30").ClearFormats()
$ws.Range("I16:I30").ClearFormats()

# Sheet3 (我的投资组合收益) - append new rows 31-45 (timestamp 202506101600)
$ws = $wb.Worksheets.Item(3)
$ws.Range("B31:B45").NumberFormat = "@"
$ws.Range("I31:I45").NumberFormat = "@"
$ws.Range("A31").Value = '范式进化投资组合'
$ws.Range("B31").Value = '000333'
$ws.Range("C31").Value = '美的集团'
$ws.Range("D31").Value = 0.4999999999999974
$ws.Range("E31").Value = 13.28200178359967
$ws.Range("F31").Value = 75.05
$ws.Range("G31").Value = 996.8142338591555
$ws.Range("H31").Value = 199362.8467718311
$ws.Range("I31").Value = '202506101600'
$ws.Range("A32").Value = '范式进化投资组合'
$ws.Range("B32").Value = '000725'
$ws.Range("C32").Value = '京东方A'
$ws.Range("D32").Value = 2.499999999999987
$ws.Range("E32").Value = 1277.966966486097
$ws.Range("F32").Value = 3.9
$ws.Range("G32").Value = 4984.071169295778
$ws.Range("H32").Value = 199362.8467718311
$ws.Range("I32").Value = '202506101600'
$ws.Range("A33").Value = '范式进化投资组合'
$ws.Range("B33").Value = '159781'
$ws.Range("C33").Value = '科创创业ETF'
$ws.Range("D33").Value = 2.499999999999987
$ws.Range("E33").Value = 9403.907866595806
$ws.Range("F33").Value = 0.53
$ws.Range("G33").Value = 4984.071169295778
$ws.Range("H33").Value = 199362.8467718311
$ws.Range("I33").Value = '202506101600'
$ws.Range("A34").Value = '范式进化投资组合'
$ws.Range("B34").Value = '510050'
$ws.Range("C34").Value = '上证50ETF'
$ws.Range("D34").Value = 2.499999999999987
$ws.Range("E34").Value = 1819.00407638532
$ws.Range("F34").Value = 2.74
$ws.Range("G34").Value = 4984.071169295778
$ws.Range("H34").Value = 199362.8467718311
$ws.Range("I34").Value = '202506101600'
$ws.Range("A35").Value = '范式进化投资组合'
$ws.Range("B35").Value = '510300'
$ws.Range("C35").Value = '沪深300ETF'
$ws.Range("D35").Value = 2.499999999999987
$ws.Range("E35").Value = 1255.433543903218
$ws.Range("F35").Value = 3.97
$ws.Range("G35").Value = 4984.071169295778
$ws.Range("H35").Value = 199362.8467718311
$ws.Range("I35").Value = '202506101600'
$ws.Range("A36").Value = '范式进化投资组合'
$ws.Range("B36").Value = '513100'
$ws.Range("C36").Value = '纳指ETF'
$ws.Range("D36").Value = 0.4999999999999974
$ws.Range("E36").Value = 634.9135247510544
$ws.Range("F36").Value = 1.57
$ws.Range("G36").Value = 996.8142338591555
$ws.Range("H36").Value = 199362.8467718311
$ws.Range("I36").Value = '202506101600'
$ws.Range("A37").Value = '范式进化投资组合'
$ws.Range("B37").Value = '513290'
$ws.Range("C37").Value = '纳指生物科技ETF'
$ws.Range("D37").Value = 0.4999999999999974
$ws.Range("E37").Value = 890.0127088028173
$ws.Range("F37").Value = 1.12
$ws.Range("G37").Value = 996.8142338591555
$ws.Range("H37").Value = 199362.8467718311
$ws.Range("I37").Value = '202506101600'
$ws.Range("A38").Value = '范式进化投资组合'
$ws.Range("B38").Value = '518880'
$ws.Range("C38").Value = '黄金ETF'
$ws.Range("D38").Value = 0.4999999999999974
$ws.Range("E38").Value = 134.5228385774839
$ws.Range("F38").Value = 7.41
$ws.Range("G38").Value = 996.8142338591555
$ws.Range("H38").Value = 199362.8467718311
$ws.Range("I38").Value = '202506101600'
$ws.Range("A39").Value = '范式进化投资组合'
$ws.Range("B39").Value = '600085'
$ws.Range("C39").Value = '同仁堂'
$ws.Range("D39").Value = 0.4999999999999974
$ws.Range("E39").Value = 26.4056750691167
$ws.Range("F39").Value = 37.75
$ws.Range("G39").Value = 996.8142338591555
$ws.Range("H39").Value = 199362.8467718311
$ws.Range("I39").Value = '202506101600'
$ws.Range("A40").Value = '范式进化投资组合'
$ws.Range("B40").Value = '600900'
$ws.Range("C40").Value = '长江电力'
$ws.Range("D40").Value = 0.4999999999999974
$ws.Range("E40").Value = 33.19394718145706
$ws.Range("F40").Value = 30.03
$ws.Range("G40").Value = 996.8142338591554
$ws.Range("H40").Value = 199362.8467718311
$ws.Range("I40").Value = '202506101600'
$ws.Range("A41").Value = '范式进化投资组合'
$ws.Range("B41").Value = '600989'
$ws.Range("C41").Value = '宝丰能源'
$ws.Range("D41").Value = 2.499999999999987
$ws.Range("E41").Value = 309.7620366249706
$ws.Range("F41").Value = 16.09
$ws.Range("G41").Value = 4984.071169295778
$ws.Range("H41").Value = 199362.8467718311
$ws.Range("I41").Value = '202506101600'
$ws.Range("A42").Value = '范式进化投资组合'
$ws.Range("B42").Value = '603119'
$ws.Range("C42").Value = '浙江荣泰'
$ws.Range("D42").Value = 0.4999999999999974
$ws.Range("E42").Value = 25.09602804277834
$ws.Range("F42").Value = 39.72
$ws.Range("G42").Value = 996.8142338591555
$ws.Range("H42").Value = 199362.8467718311
$ws.Range("I42").Value = '202506101600'
$ws.Range("A43").Value = '范式进化投资组合'
$ws.Range("B43").Value = 'HK02899'
$ws.Range("C43").Value = '紫金矿业'
$ws.Range("D43").Value = 0.4999999999999974
$ws.Range("E43").Value = 54.17468662278019
$ws.Range("F43").Value = 18.4
$ws.Range("G43").Value = 996.8142338591555
$ws.Range("H43").Value = 199362.8467718311
$ws.Range("I43").Value = '202506101600'
$ws.Range("A44").Value = '范式进化投资组合'
$ws.Range("B44").Value = 'HK06881'
$ws.Range("C44").Value = '中国银河'
$ws.Range("D44").Value = 0.4999999999999974
$ws.Range("E44").Value = 122.7603736279748
$ws.Range("F44").Value = 8.12
$ws.Range("G44").Value = 996.8142338591555
$ws.Range("H44").Value = 199362.8467718311
$ws.Range("I44").Value = '202506101600'
$ws.Range("A45").Value = '范式进化投资组合'
$ws.Range("B45").Value = '100000'
$ws.Range("C45").Value = '现金'
$ws.Range("D45").Value = 32.99999999999984
$ws.Range("E45").Value = 65789.73943470427
$ws.Range("F45").Value = 1.0
$ws.Range("G45").Value = 65789.73943470427
$ws.Range("H45").Value = 199362.8467718311
$ws.Range("I45").Value = '202506101600'
$ws.Range("B31:B45").ClearFormats()
$ws.Range("I31:I45").ClearFormats()
